$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = '66.135.20'
$ws.Range("E2").Value = '  -0.63%  '
$ws.Range("D3").Value = '3.226.07'
$ws.Range("E3").Value = '  -3.61%  '
$ws.Range("E4").Value = '  +0.13%  '
$c = $ws.Range("D5")
$c.NumberFormat = '@'
$c.Value = '574.96'
$c.Style = 'Normal'
$ws.Range("E5").Value = '  -1.94%  '
$c = $ws.Range("D6")
$c.NumberFormat = '@'
$c.Value = '168.99'
$c.Style = 'Normal'
$ws.Range("E6").Value = '  -9.43%  '
$ws.Range("E7").Value = '  +0.15%  '
$c = $ws.Range("D8")
$c.NumberFormat = '@'
$c.Value = '0.572'
$c.Style = 'Normal'
$ws.Range("E8").Value = '  -1.02%  '
$ws.Range("D9").Value = '3.220.77'
$ws.Range("E9").Value = '  -3.61%  '
$c = $ws.Range("D10")
$c.NumberFormat = '@'
$c.Value = '0.169'
$c.Style = 'Normal'
$ws.Range("E10").Value = '  -7.72%  '
$c = $ws.Range("D11")
$c.NumberFormat = '@'
$c.Value = '0.563'
$c.Style = 'Normal'
$ws.Range("E11").Value = '  -3.72%  '
$c = $ws.Range("D12")
$c.NumberFormat = '@'
$c.Value = '44.16'
$c.Style = 'Normal'
$ws.Range("E12").Value = '  -6.49%  '
$c = $ws.Range("D13")
$c.NumberFormat = '@'
$c.Value = '0.0000266'
$c.Style = 'Normal'
$ws.Range("E13").Value = '  -2.29%  '
$c = $ws.Range("D14")
$c.NumberFormat = '@'
$c.Value = '667.89'
$c.Style = 'Normal'
$ws.Range("E14").Value = '  +2.68%  '
$ws.Range("D15").Value = '3.755.97'
$ws.Range("E15").Value = '  -3.24%  '
$c = $ws.Range("D16")
$c.NumberFormat = '@'
$c.Value = '8.14'
$c.Style = 'Normal'
$ws.Range("E16").Value = '  -4.40%  '
$ws.Range("D17").Value = '66.170.34'
$ws.Range("E17").Value = '  -0.60%  '
$ws.Range("E18").Value = '  -0.31%  '
$ws.Range("D19").Value = '3.234.98'
$ws.Range("E19").Value = '  -3.33%  '
$c = $ws.Range("D20")
$c.NumberFormat = '@'
$c.Value = '16.98'
$c.Style = 'Normal'
$ws.Range("E20").Value = '  -5.36%  '
$c = $ws.Range("D21")
$c.NumberFormat = '@'
$c.Value = '10.55'
$c.Style = 'Normal'
$ws.Range("E21").Value = '  -5.43%  '
$c = $ws.Range("D22")
$c.NumberFormat = '@'
$c.Value = '0.864'
$c.Style = 'Normal'
$ws.Range("E22").Value = '  -4.07%  '
$c = $ws.Range("D23")
$c.NumberFormat = '@'
$c.Value = '16.71'
$c.Style = 'Normal'
$ws.Range("E23").Value = '  -5.87%  '
$c = $ws.Range("D24")
$c.NumberFormat = '@'
$c.Value = '5.19'
$c.Style = 'Normal'
$ws.Range("E24").Value = '  +1.98%  '
$c = $ws.Range("D25")
$c.NumberFormat = '@'
$c.Value = '95.57'
$c.Style = 'Normal'
$ws.Range("E25").Value = '  -4.65%  '
$c = $ws.Range("D26")
$c.NumberFormat = '@'
$c.Value = '3.79'
$c.Style = 'Normal'
$ws.Range("E26").Value = '  -5.20%  '
$c = $ws.Range("D27")
$c.NumberFormat = '@'
$c.Value = '2.60'
$c.Style = 'Normal'
$ws.Range("E27").Value = '  -7.31%  '
$c = $ws.Range("D28")
$c.NumberFormat = '@'
$c.Value = '8.89'
$c.Style = 'Normal'
$ws.Range("E28").Value = '  -8.08%  '
$c = $ws.Range("D29")
$c.NumberFormat = '@'
$c.Value = '31.79'
$c.Style = 'Normal'
$ws.Range("E29").Value = '  -1.07%  '
$c = $ws.Range("D30")
$c.NumberFormat = '@'
$c.Value = '8.15'
$c.Style = 'Normal'
$ws.Range("E30").Value = '  -4.77%  '
$ws.Range("E31").Value = '  -5.19%  '
$c = $ws.Range("D32")
$c.NumberFormat = '@'
$c.Value = '561.70'
$c.Style = 'Normal'
$ws.Range("E32").Value = '  -8.01%  '
$c = $ws.Range("D33")
$c.NumberFormat = '@'
$c.Value = '10.72'
$c.Style = 'Normal'
$ws.Range("E33").Value = '  -3.76%  '
$ws.Range("D34").Value = '3.775.86'
$ws.Range("E34").Value = '  -2.87%  '
$c = $ws.Range("D35")
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range("E35").Value = '  -0.04%  '
$ws.Range("E36").Value = '  -4.68%  '
$ws.Range("B37").Value = 'dogwifhat'
$ws.Range("C37").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$c = $ws.Range("D37")
$c.NumberFormat = '@'
$c.Value = '3.29'
$c.Style = 'Normal'
$ws.Range("E37").Value = '  -16.91%  '
$ws.Range("B38").Value = 'OKB'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Range("D38")
$c.NumberFormat = '@'
$c.Value = '54.73'
$c.Style = 'Normal'
$ws.Range("E38").Value = '  -2.39%  '
$c = $ws.Range("D39")
$c.NumberFormat = '@'
$c.Value = '0.127'
$c.Style = 'Normal'
$ws.Range("E39").Value = '  -2.57%  '
$ws.Range("B40").Value = 'ApeXProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$c = $ws.Range("D40")
$c.NumberFormat = '@'
$c.Value = '3.30'
$c.Style = 'Normal'
$ws.Range("E40").Value = '  -2.00%  '
$ws.Range("B41").Value = 'Fetch.AI'
$ws.Range("C41").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c = $ws.Range("D41")
$c.NumberFormat = '@'
$c.Value = '2.53'
$c.Style = 'Normal'
$ws.Range("E41").Value = '  -7.87%  '
$ws.Range("B42").Value = 'InjectiveProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Range("D42")
$c.NumberFormat = '@'
$c.Value = '31.22'
$c.Style = 'Normal'
$ws.Range("E42").Value = '  -6.62%  '
$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range("D43")
$c.NumberFormat = '@'
$c.Value = '2.94'
$c.Style = 'Normal'
$ws.Range("E43").Value = '  -8.54%  '
$ws.Range("B44").Value = 'PEPE'
$ws.Range("C44").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D44").Value = '0.0₃0640'
$ws.Range("E44").Value = '  -8.78%  '
$c = $ws.Range("D45")
$c.NumberFormat = '@'
$c.Value = '0.318'
$c.Style = 'Normal'
$ws.Range("E45").Value = '  -7.25%  '
$c = $ws.Range("D46")
$c.NumberFormat = '@'
$c.Value = '0.0391'
$c.Style = 'Normal'
$ws.Range("E46").Value = '  -6.56%  '
$c = $ws.Range("D47")
$c.NumberFormat = '@'
$c.Value = '1.01'
$c.Style = 'Normal'
$ws.Range("E47").Value = '  +0.28%  '
$ws.Range("E48").Value = '  -2.76%  '
$ws.Range("E49").Value = '  -2.99%  '
$ws.Range("E50").Value = '  -5.81%  '
$c = $ws.Range("D51")
$c.NumberFormat = '@'
$c.Value = '126.89'
$c.Style = 'Normal'
$ws.Range("E51").Value = '  -1.93%  '
